$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Be able to create account. ... connected to the databases" paragraph:
#    split the runs around "account" and "databases" and wrap each in
#    w:proofErr gramStart/gramEnd markers (as Word's grammar checker would).
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Be able to create account*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' +
    '<w:p w14:paraId="1E793020" w14:textId="25D790C5" w:rsidR="00D20ACA" w:rsidRDefault="00D20ACA" w:rsidP="00D20ACA">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Be able to create </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>account</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:t xml:space="preserve">. Each account will be able to post items to sell them or add items to a cart to buy them. </w:t></w:r>' +
      '<w:r w:rsidR="00586314"><w:t>Application uses real life cryptocurrencies for transactions</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> -&gt; all these actions are connected to the </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>databases</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

    $target.Range.InsertXML($fragment)
}

# ---------------------------------------------------------------------------
# 2) Locate the paragraph that holds the first (MySQL) hyperlink so we can
#    append the two new paragraphs right after it, before the document's
#    trailing empty paragraph.
# ---------------------------------------------------------------------------
$mysqlHyperlinkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*3vsC05rxZ8c*") {
        $mysqlHyperlinkPara = $p
        break
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# ---------------------------------------------------------------------------
# 3) New bold heading paragraph: "Django Tutorial for Beginners ..."
#    Inserted immediately before the document's trailing empty paragraph.
# ---------------------------------------------------------------------------
$headingFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
  '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Django Tutorial for Beginners' + [char]0x2013 + ' Build Powerful Backends </w:t></w:r>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$lastPara.Range.InsertXML($headingFragment)

# ---------------------------------------------------------------------------
# 4) New hyperlink paragraph pointing at the Django tutorial video. Re-use
#    the (now shifted) trailing empty paragraph for it, exactly like Word
#    would when the cursor sits on the last empty line and a link is
#    inserted, then restore a fresh trailing empty paragraph afterwards.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$null = $d.Hyperlinks.Add($lastPara.Range, "https://www.youtube.com/watch?v=rHux0gMZ3Eg&t=160s")

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Output "done"
